$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44211
$ws.Range("M2").Value = 200
$ws.Range("N2").Value = 3000
$ws.Range("O2").Value = 3500
$ws.Range("P2").Value = 3250
$ws.Range("R2").Value = 'Provincia de Curicó'
$ws.Range("S2").Value = 1625
$ws.Range("D3").Value = 44204
$ws.Range("M3").Value = 150
$ws.Range("N3").Value = 4000
$ws.Range("O3").Value = 4000
$ws.Range("P3").Value = 4000
$ws.Range("S3").Value = 2000
$ws.Range("D4").Value = 44204
$ws.Range("M4").Value = 250
$ws.Range("D5").Value = 44188
$ws.Range("L5").Value = 'Primera'
$ws.Range("M5").Value = 300
$ws.Range("N5").Value = 4000
$ws.Range("O5").Value = 4000
$ws.Range("P5").Value = 4000
$ws.Range("R5").Value = 'Provincia de Curicó'
$ws.Range("S5").Value = 2000
$ws.Range("D6").Value = 44188
$ws.Range("M6").Value = 500
$ws.Range("N6").Value = 4000
$ws.Range("O6").Value = 4000
$ws.Range("P6").Value = 4000
$ws.Range("R6").Value = 'Provincia de Linares'
$ws.Range("S6").Value = 2000
$ws.Range("D7").Value = 44225
$ws.Range("M7").Value = 150
$ws.Range("N7").Value = 4000
$ws.Range("O7").Value = 4000
$ws.Range("P7").Value = 4000
$ws.Range("R7").Value = 'Provincia de Curicó'
$ws.Range("S7").Value = 2000
$ws.Range("D8").Value = 44225
$ws.Range("M8").Value = 200
$ws.Range("N8").Value = 4000
$ws.Range("P8").Value = 4000
$ws.Range("S8").Value = 2000
$ws.Range("D9").Value = 44193
$ws.Range("M9").Value = 200
$ws.Range("N9").Value = 3000
$ws.Range("O9").Value = 3000
$ws.Range("P9").Value = 3000
$ws.Range("R9").Value = 'Provincia de Linares'
$ws.Range("S9").Value = 1500
$ws.Range("D10").Value = 44201
$ws.Range("M10").Value = 200
$ws.Range("R10").Value = 'Provincia de Linares'
$ws.Range("D11").Value = 44202
$ws.Range("M11").Value = 200
$ws.Range("N11").Value = 4000
$ws.Range("P11").Value = 4000
$ws.Range("R11").Value = 'Provincia de Linares'
$ws.Range("S11").Value = 2000
$ws.Range("D12").Value = 44169
$ws.Range("M12").Value = 200
$ws.Range("N12").Value = 5000
$ws.Range("O12").Value = 5000
$ws.Range("P12").Value = 5000
$ws.Range("R12").Value = 'Provincia de Linares'
$ws.Range("S12").Value = 2500
$ws.Range("D13").Value = 44216
$ws.Range("R13").Value = 'Provincia de Curicó'
$ws.Range("D14").Value = 44216
$ws.Range("M14").Value = 400
$ws.Range("D15").Value = 44260
$ws.Range("M15").Value = 75
$ws.Range("D16").Value = 44189
$ws.Range("M16").Value = 300
$ws.Range("N16").Value = 3000
$ws.Range("O16").Value = 3000
$ws.Range("P16").Value = 3000
$ws.Range("R16").Value = 'Provincia de Curicó'
$ws.Range("S16").Value = 1500
$ws.Range("D17").Value = 44189
$ws.Range("M17").Value = 250
$ws.Range("N17").Value = 3000
$ws.Range("O17").Value = 3000
$ws.Range("P17").Value = 3000
$ws.Range("R17").Value = 'Provincia de Linares'
$ws.Range("S17").Value = 1500
$ws.Range("D18").Value = 44181
$ws.Range("M18").Value = 140
$ws.Range("O18").Value = 4500
$ws.Range("P18").Value = 4250
$ws.Range("S18").Value = 2125
$ws.Range("D19").Value = 44203
$ws.Range("M19").Value = 350
$ws.Range("D20").Value = 44239
$ws.Range("M20").Value = 350
$ws.Range("N20").Value = 3500
$ws.Range("P20").Value = 3750
$ws.Range("S20").Value = 1875
$ws.Range("D21").Value = 44187
$ws.Range("M21").Value = 100
$ws.Range("N21").Value = 3400
$ws.Range("O21").Value = 3400
$ws.Range("P21").Value = 3400
$ws.Range("R21").Value = 'Provincia de Curicó'
$ws.Range("S21").Value = 1700
$ws.Range("D22").Value = 44187
$ws.Range("M22").Value = 200
$ws.Range("R22").Value = 'Provincia de Linares'
$ws.Range("D23").Value = 44187
$ws.Range("L23").Value = 'Segunda'
$ws.Range("M23").Value = 50
$ws.Range("N23").Value = 3000
$ws.Range("O23").Value = 3000
$ws.Range("P23").Value = 3000
$ws.Range("R23").Value = 'Provincia de Linares'
$ws.Range("S23").Value = 1500
$ws.Range("D24").Value = 44195
$ws.Range("M24").Value = 300
$ws.Range("N24").Value = 3000
$ws.Range("O24").Value = 3000
$ws.Range("P24").Value = 3000
$ws.Range("S24").Value = 1500
$ws.Range("D25").Value = 44209
$ws.Range("M25").Value = 170
$ws.Range("N25").Value = 3000
$ws.Range("P25").Value = 3500
$ws.Range("S25").Value = 1750
$ws.Range("D26").Value = 44250
$ws.Range("M26").Value = 100
$ws.Range("N26").Value = 4000
$ws.Range("O26").Value = 4000
$ws.Range("P26").Value = 4000
$ws.Range("R26").Value = 'Provincia de Curicó'
$ws.Range("S26").Value = 2000
$ws.Range("D27").Value = 44221
$ws.Range("M27").Value = 150
$ws.Range("D28").Value = 44221
$ws.Range("M28").Value = 200
$ws.Range("D29").Value = 44186
$ws.Range("M29").Value = 200
$ws.Range("R29").Value = 'Provincia de Linares'
$ws.Range("D30").Value = 44252
$ws.Range("M30").Value = 75
$ws.Range("N30").Value = 4000
$ws.Range("P30").Value = 4000
$ws.Range("R30").Value = 'Provincia de Curicó'
$ws.Range("S30").Value = 2000
$ws.Range("D31").Value = 44224
$ws.Range("M31").Value = 250
$ws.Range("N31").Value = 4000
$ws.Range("O31").Value = 4000
$ws.Range("P31").Value = 4000
$ws.Range("R31").Value = 'Provincia de Curicó'
$ws.Range("S31").Value = 2000
$ws.Range("D32").Value = 44224
$ws.Range("M32").Value = 300
$ws.Range("D33").Value = 44217
$ws.Range("M33").Value = 250
$ws.Range("D34").Value = 44217
$ws.Range("M34").Value = 300
$ws.Range("R34").Value = 'Provincia de Linares'
$ws.Range("D35").Value = 44196
$ws.Range("M35").Value = 150
$ws.Range("D36").Value = 44215
$ws.Range("M36").Value = 750
$ws.Range("O36").Value = 4000
$ws.Range("P36").Value = 4000
$ws.Range("R36").Value = 'Provincia de Curicó'
$ws.Range("S36").Value = 2000
$ws.Range("D37").Value = 44194
$ws.Range("M37").Value = 250
$ws.Range("N37").Value = 4000
$ws.Range("O37").Value = 4000
$ws.Range("P37").Value = 4000
$ws.Range("S37").Value = 2000
$ws.Range("D38").Value = 44222
$ws.Range("M38").Value = 250
$ws.Range("R38").Value = 'Provincia de Curicó'
$ws.Range("D39").Value = 44222
$ws.Range("R39").Value = 'Provincia de Linares'
$ws.Range("D40").Value = 44210
$ws.Range("M40").Value = 400
$ws.Range("N40").Value = 3000
$ws.Range("P40").Value = 3500
$ws.Range("S40").Value = 1750
$ws.Range("D42").Value = 44176
$ws.Range("M42").Value = 100
$ws.Range("D43").Value = 44257
$ws.Range("M43").Value = 100
$ws.Range("N43").Value = 4000
$ws.Range("O43").Value = 4000
$ws.Range("P43").Value = 4000
$ws.Range("S43").Value = 2000
$ws.Range("D44").Value = 44251
$ws.Range("M44").Value = 125
$ws.Range("D45").Value = 44175
$ws.Range("M45").Value = 250
